# Generate Report for Handback
#
# The handback process detected that the file name returned by the
# translation vendor didn't match the file name that was originally
# handed off, so the status for both locales flips from "Ready for
# handoff" to "Handback transform failed" and an explanatory message is
# recorded in the (previously empty) "Error Detail" column for the
# second (73a1a76b...) row of each locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column ("Ready for handoff" -> "Handback transform failed") for
# the 73a1a76b... row on every sheet that surfaces it.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhcn.Range("C3").Value     = "Handback transform failed"
$dede.Range("C3").Value     = "Handback transform failed"

# Error Detail column (L) explaining the handback/handoff file name
# mismatch, one message per locale.
$zhcn.Range("L3").Value = "Handback file name: ysphktuv.r0k is different with handoff file name: 73a1a76b-1f3d-445b-8126-9a8fc337cd24.cd0db91349a944fa5ee9fa4b6def4b3b1f19bc72.zh-cn."
$dede.Range("L3").Value = "Handback file name: ysphktuv.r0k is different with handoff file name: 73a1a76b-1f3d-445b-8126-9a8fc337cd24.cd0db91349a944fa5ee9fa4b6def4b3b1f19bc72.de-de."
